$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23: "Be the only one" blog post update
$ws.Range("D23").Value = "안녕하세요! `n`n카사바 잎 질병 분류 대회(Cassava Leaf Disease Classification Competition)가 끝나고 개인적"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2733"

# Row 29: promedius.ai blog post update
$ws.Range("D29").Value = "TorchIO를 이용한 3D Segmentation"
$ws.Range("E29").Value = "https://blog.promedius.ai/torchioreul-iyonghan-3d-segmentation/"

# Row 37: DSBA seminar update
$ws.Range("D37").Value = "[Paper Review] BatchEnsemble: An Alternative Approach to Efficient Ensemble and Lifelong Learning"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1445&mod=document&pageid=1"

# Row 39: a292run.tistory.com blog post update
$ws.Range("D39").Value = "Deep Dive Analysis of Missing Values in Dataset"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Deep-Dive-Analysis-of-Missing-Values-in-Dataset-1"

# Row 51: bskyvision.com blog post update
$ws.Range("D51").Value = '[독후감] 거리 위 청소년들의 아빠, 이요셉 목사의 <지금 가고 있어>'
$ws.Range("E51").Value = "https://bskyvision.com/983"
